# Added 13th Test Case in feature file, steps and Home Page
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXHome")

# Populate the newly added 13th test case data in rows 11 and 12
$ws.Range("B11").Value = "Spot"
$ws.Range("C11").Value = "Trader01@Tinyex"
$ws.Range("D11").Value = "ETH"
$ws.Range("E11").Value = "USDT"
$ws.Range("J11").Value = 1

$ws.Range("B12").Value = "Spot"
$ws.Range("C12").Value = "Trader01@Tinyex"
$ws.Range("D12").Value = "ETH"
$ws.Range("E12").Value = "USDT"
$ws.Range("J12").Value = 1

# Update the active selection to match the edited range
$ws.Activate()
$ws.Range("B11:E11").Select()
